$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "vehicle"

$ws.Range("D5").Value = "smallint"
$ws.Range("E5").Value = ""
$ws.Range("D37").Value = "Varchar"
$ws.Range("D60").Value = "TblLocation"

$null = $ws.Range("D60:F60").Select()
